$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this shifts existing rows 28..130 down to 29..131,
# preserving all of their values and formatting (including the date style on column D).
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly price record.
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C28").Value = "Los Lagos"
$ws.Range("D28").Value = 44481
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 100112032
$ws.Range("G28").Value = "Zapallo italiano"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 15000
$ws.Range("N28").Value = "$/caja 50 unidades"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 300
$ws.Range("Q28").Value = 50
$ws.Range("R28").Value = "Hortaliza"
